$d = $word.ActiveDocument

# Locate the relevant paragraphs by their text rather than a hard-coded
# index, so the script is resilient to any unrelated paragraph shuffling.
$greenIdx = 0
$instrIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.TrimEnd()
    if ($t -eq "Green Tomatoes\n") {
        $greenIdx = $i
    }
    if ($t -eq "Instructions") {
        $instrIdx = $i
    }
}

# --- 1. "Green Tomatoes" paragraph + following blank "\n" paragraph ---
# Before:
#   Para: ind left=720, runs: "Green " + "Tomatoes\n"
#   Para: ind firstLine=720, run: "\n"
# After:
#   Para: (no ind), runs: "                " (16 spaces, no rPr) + "Green Tomatoes\n"
#   Para: ind firstLine=720, run: "\n"   (same as the old second paragraph)
$pGreen = $d.Paragraphs($greenIdx)
$pBlank = $d.Paragraphs($greenIdx + 1)
$combined = $d.Range($pGreen.Range.Start, $pBlank.Range.End)

$xmlGreen = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:rPr><w:color w:val="414141"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">                </w:t></w:r><w:r><w:rPr><w:color w:val="414141"/></w:rPr><w:t>Green Tomatoes\n</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:ind w:firstLine="720"/><w:rPr><w:color w:val="414141"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="414141"/></w:rPr><w:t>\n</w:t></w:r></w:p>
'@

$combined.InsertXML($xmlGreen)

# --- 2. "Instructions" paragraph gains a trailing "\n" run ---
# Before: run: "Instructions"
# After:  runs: "Instructions" + "\n"
$pInstr = $d.Paragraphs($instrIdx)

$xmlInstr = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:rPr><w:color w:val="414141"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="414141"/></w:rPr><w:t>Instructions</w:t></w:r><w:r><w:rPr><w:color w:val="414141"/></w:rPr><w:t>\n</w:t></w:r></w:p>
'@

$pInstr.Range.InsertXML($xmlInstr)
